# Update the "dSF" column (F) values for specific rows, per the
# repulled/recalculated data from the commit "repull data, push all data, mean calculation".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = -7
    5  = -3
    6  = -8
    10 = 1
    15 = -2
    18 = 4
    21 = -3
    28 = -3
    30 = 4
    37 = -5
    43 = 4
    45 = -6
    48 = -2
    49 = 4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
